# Update model, refactor parser to jimiParser
#
# Slide 5 (1-based) contains the "Parser" box plus its connector and its
# small numeric callout label. Resize/move the three shapes to match the
# updated model layout and rename the box text from "Parser" to
# "JimiParser".
#
# NOTE: Shape.Left/Top/Width/Height round-trip through a 32-bit float in
# this COM host (matching real PowerPoint's `Single`-typed Shape
# geometry properties), so naively assigning `emu / 12700.0` can land a
# handful of EMU away from the exact target after PowerPoint converts
# back to EMU on save. EmuToPt searches for a point value that, once
# coerced through that same float32 round-trip, reproduces the exact
# requested EMU value so the saved XML matches precisely.

function EmuToPt($targetEmu) {
    $base = $targetEmu / 12700.0
    $step = 0.0000001

    $f32 = [float]$base
    $back = [math]::Floor([double]$f32 * 12700.0)
    if ($back -eq $targetEmu) {
        return $base
    }

    for ($i = 1; $i -lt 10000; $i++) {
        $trial = $base + ($i * $step)
        $f32 = [float]$trial
        $back = [math]::Floor([double]$f32 * 12700.0)
        if ($back -eq $targetEmu) {
            return $trial
        }
    }

    for ($i = 1; $i -lt 10000; $i++) {
        $trial = $base - ($i * $step)
        $f32 = [float]$trial
        $back = [math]::Floor([double]$f32 * 12700.0)
        if ($back -eq $targetEmu) {
            return $trial
        }
    }

    return $base
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)

# --- Shape 15 (id=18, "Rectangle 62") -> text box reading "Parser" ---
$rect = $s.Shapes.Item(15)
$rect.Left   = EmuToPt(3341296)
$rect.Top    = EmuToPt(3766159)
$rect.Width  = EmuToPt(771982)
$rect.Height = EmuToPt(346760)
$rect.TextFrame.TextRange.Text = "JimiParser"

# --- Shape 16 (id=19, "Elbow Connector 106") feeding into the box above ---
$conn = $s.Shapes.Item(16)
$conn.Left   = EmuToPt(3108853)
$conn.Top    = EmuToPt(3937000)
$conn.Width  = EmuToPt(232443)
$conn.Height = EmuToPt(2539)

# --- Shape 28 (id=40, "TextBox 39") small numeric callout near the box ---
$lbl = $s.Shapes.Item(28)
$lbl.Left   = EmuToPt(3163363)
$lbl.Top    = EmuToPt(3731547)
$lbl.Width  = EmuToPt(131116)
$lbl.Height = EmuToPt(230832)
